$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 2 (shifts existing rows 2-6 down to 3-7)
$ws.Rows.Item(2).Insert()

# The inserted row copies formatting from the row above (the bold/bordered
# header row); data rows in this sheet carry no explicit style, so strip it.
$ws.Rows.Item(2).ClearFormats()

# Populate the newly inserted row 2 with the new match data
$ws.Cells.Item(2, 1).Value = 'EouYkCd3'
$ws.Cells.Item(2, 2).Value = '28/11/2024'
$ws.Cells.Item(2, 3).Value = '16:00'
$ws.Cells.Item(2, 4).Value = 'BOLIVIA - DIVISION PROFESIONAL'
$ws.Cells.Item(2, 5).Value = 'Santa Cruz'
$ws.Cells.Item(2, 6).Value = 'Independiente'
$ws.Cells.Item(2, 7).Value = 2
$ws.Cells.Item(2, 8).Value = 3.6
$ws.Cells.Item(2, 9).Value = 3.5
$ws.Cells.Item(2, 10).Value = 2.63
$ws.Cells.Item(2, 11).Value = 2.2
$ws.Cells.Item(2, 12).Value = 4
$ws.Cells.Item(2, 13).Value = 1.04
$ws.Cells.Item(2, 14).Value = 12
$ws.Cells.Item(2, 15).Value = 1.25
$ws.Cells.Item(2, 16).Value = 3.75
$ws.Cells.Item(2, 17).Value = 1.83
$ws.Cells.Item(2, 18).Value = 1.98
$ws.Cells.Item(2, 19).Value = 1.36
$ws.Cells.Item(2, 20).Value = 3
$ws.Cells.Item(2, 21).Value = 1.73
$ws.Cells.Item(2, 22).Value = 2
$ws.Cells.Item(2, 23).Value = 8
$ws.Cells.Item(2, 24).Value = 10
$ws.Cells.Item(2, 25).Value = 9
$ws.Cells.Item(2, 26).Value = 17
$ws.Cells.Item(2, 27).Value = 15
$ws.Cells.Item(2, 28).Value = 23
$ws.Cells.Item(2, 29).Value = 12
$ws.Cells.Item(2, 30).Value = 7
$ws.Cells.Item(2, 31).Value = 13
$ws.Cells.Item(2, 32).Value = 41
$ws.Cells.Item(2, 33).Value = 201
$ws.Cells.Item(2, 34).Value = 12
$ws.Cells.Item(2, 35).Value = 19
$ws.Cells.Item(2, 36).Value = 12
$ws.Cells.Item(2, 37).Value = 41
$ws.Cells.Item(2, 38).Value = 29
$ws.Cells.Item(2, 39).Value = 34
$ws.Cells.Item(2, 40).Value = 4
$ws.Cells.Item(2, 41).Value = 11
$ws.Cells.Item(2, 42).Value = 21
$ws.Cells.Item(2, 43).Value = 34
$ws.Cells.Item(2, 44).Value = 51
$ws.Cells.Item(2, 45).Value = 126
$ws.Cells.Item(2, 46).Value = 3
$ws.Cells.Item(2, 47).Value = 7.5
$ws.Cells.Item(2, 48).Value = 51
$ws.Cells.Item(2, 49).Value = 5.5
$ws.Cells.Item(2, 50).Value = 19
$ws.Cells.Item(2, 51).Value = 26
$ws.Cells.Item(2, 52).Value = 67
$ws.Cells.Item(2, 53).Value = 81
$ws.Cells.Item(2, 54).Value = 151
$ws.Cells.Item(2, 55).Value = ''
$ws.Cells.Item(2, 56).Value = ''
